$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Cell='D2'; Value='27.045.98'},
  @{Cell='E2'; Value='  +0.62%  '},
  @{Cell='D3'; Value='1.682.65'},
  @{Cell='E3'; Value='  +0.90%  '},
  @{Cell='E4'; Value='  +0.03%  '},
  @{Cell='D5'; Value='216.13'},
  @{Cell='E5'; Value='  +0.26%  '},
  @{Cell='E6'; Value='  -2.11%  '},
  @{Cell='E7'; Value='  -0.04%  '},
  @{Cell='E8'; Value='  +0.32%  '},
  @{Cell='D9'; Value='21.51'},
  @{Cell='E9'; Value='  +6.40%  '},
  @{Cell='D10'; Value='0.0624'},
  @{Cell='E10'; Value='  +0.82%  '},
  @{Cell='D11'; Value='0.0890'},
  @{Cell='E11'; Value='  -0.66%  '},
  @{Cell='D12'; Value='1.921.31'},
  @{Cell='E12'; Value='  +0.99%  '},
  @{Cell='D13'; Value='1.694.39'},
  @{Cell='E13'; Value='  +1.65%  '},
  @{Cell='E14'; Value='  +0.61%  '},
  @{Cell='E15'; Value='  +1.93%  '},
  @{Cell='E16'; Value='  +0.65%  '},
  @{Cell='D17'; Value='27.058.65'},
  @{Cell='E17'; Value='  +0.59%  '},
  @{Cell='D18'; Value='8.18'},
  @{Cell='E18'; Value='  +4.97%  '},
  @{Cell='D19'; Value='236.52'},
  @{Cell='E19'; Value='  +2.09%  '},
  @{Cell='D20'; Value='0.0₃0739'},
  @{Cell='E20'; Value='  +0.80%  '},
  @{Cell='E21'; Value='  +0.03%  '},
  @{Cell='D22'; Value='4.48'},
  @{Cell='E22'; Value='  +0.39%  '},
  @{Cell='D23'; Value='9.29'},
  @{Cell='E23'; Value='  +1.18%  '},
  @{Cell='E24'; Value='  -3.70%  '},
  @{Cell='D25'; Value='147.21'},
  @{Cell='E25'; Value='  +1.08%  '},
  @{Cell='E26'; Value='  +5.45%  '},
  @{Cell='E27'; Value='  +1.69%  '},
  @{Cell='E28'; Value='  -2.35%  '},
  @{Cell='E29'; Value='  -0.02%  '},
  @{Cell='E30'; Value='  +0.57%  '},
  @{Cell='E31'; Value='  +0.08%  '},
  @{Cell='D32'; Value='3.36'},
  @{Cell='E32'; Value='  +0.59%  '},
  @{Cell='D33'; Value='1.531.06'},
  @{Cell='E33'; Value='  +4.52%  '},
  @{Cell='D34'; Value='3.18'},
  @{Cell='E34'; Value='  +0.80%  '},
  @{Cell='E35'; Value='  +4.78%  '},
  @{Cell='D36'; Value='2.40'},
  @{Cell='E36'; Value='  -0.40%  '},
  @{Cell='D37'; Value='0.591'},
  @{Cell='E37'; Value='  +3.32%  '},
  @{Cell='D38'; Value='0.919'},
  @{Cell='E38'; Value='  +2.42%  '},
  @{Cell='E39'; Value='  +3.57%  '},
  @{Cell='D40'; Value='1.04'},
  @{Cell='E40'; Value='  +6.00%  '},
  @{Cell='E41'; Value='  -1.10%  '},
  @{Cell='E42'; Value='  -0.04%  '},
  @{Cell='D43'; Value='68.06'},
  @{Cell='E43'; Value='  +3.55%  '},
  @{Cell='E44'; Value='  -0.59%  '},
  @{Cell='D45'; Value='1.825.68'},
  @{Cell='E45'; Value='  +0.59%  '},
  @{Cell='E46'; Value='  +0.11%  '},
  @{Cell='D47'; Value='90.42'},
  @{Cell='D48'; Value='0.105'},
  @{Cell='E48'; Value='  +4.19%  '},
  @{Cell='E49'; Value='  +0.17%  '},
  @{Cell='B50'; Value='BabyDogeCoin'},
  @{Cell='C50'; Value='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'},
  @{Cell='D50'; Value='0.0₆0104'},
  @{Cell='E50'; Value='  -1.20%  '},
  @{Cell='B51'; Value='EnergySwap'},
  @{Cell='C51'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
  @{Cell='D51'; Value='7.90'},
  @{Cell='E51'; Value='  +4.68%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
    $ws.Range($u.Cell).Style = "Normal"
}
